$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing X9/Y9 values on the existing last row
$ws.Range("X9").Value = -0.21000099999999833
$ws.Range("Y9").Value = "Down"

# Append a new row of data (row 10)
$ws.Range("A10").Value = 42653.879479166666
$ws.Range("A10").NumberFormat = "m/d/yy h:mm"
$ws.Range("B10").Value = 13
$ws.Range("C10").Value = "Buy"
$ws.Range("D10").Value = 38
$ws.Range("E10").Value = 6667
$ws.Range("F10").Value = 337
$ws.Range("G10").Value = 63
$ws.Range("H10").Value = 35
$ws.Range("I10").Value = 92
$ws.Range("J10").Value = 7
$ws.Range("K10").Value = 9772
$ws.Range("L10").Value = 94
$ws.Range("M10").Value = 52
$ws.Range("N10").Value = 35
$ws.Range("O10").Value = 3
$ws.Range("P10").Value = "Named"
$ws.Range("Q10").Value = 47.96375473473072
$ws.Range("R10").Value = 0.49
$ws.Range("S10").Value = 0.0521
$ws.Range("S10").NumberFormat = "0.00%"
$ws.Range("T10").Value = -0.0214
$ws.Range("T10").NumberFormat = "0.00%"
$ws.Range("U10").Value = 2.25
$ws.Range("V10").Value = "N/A"
$ws.Range("W10").Value = 0
